# Menu items view module update: sets Name / Description / Price text for
# the three "Specials" rows and repositions/resizes the affected textboxes
# to match their new (possibly multi-line) content.

function EMU($v) {
    # PowerPoint COM Left/Top/Width/Height are expressed in points; the
    # host truncates points*12700 when converting back to EMU, so nudge by
    # half an EMU to land on the exact target unit.
    return ($v + 0.5) / 12700.0
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Gray used for the description lines (tx1 @ 50% luminance mod/offset).
$descGray = 8421504  # RGB(128,128,128) = 0x808080

# ---- Item 1 : Name ----------------------------------------------------
$shp = $s.Shapes.Item("TextBox 4")
$shp.TextFrame.TextRange.Delete()
$shp.TextFrame.TextRange.Text = "Club Sandwich"
$shp.TextFrame.TextRange.Font.Bold = $true
$shp.Left = EMU(3019393)
$shp.Top = EMU(1567259)
$shp.Width = EMU(5539563)
$shp.Height = EMU(369332)

# ---- Item 1 : Description ---------------------------------------------
$shp = $s.Shapes.Item("TextBox 13")
$shp.TextFrame.TextRange.Delete()
$shp.TextFrame.TextRange.Text = "chicken breast, bacon, tomato, lettuce, cheese, avocado & aioli on turkish with chips"
$shp.TextFrame.TextRange.Font.Size = 16
$shp.TextFrame.TextRange.Font.Color.RGB = $descGray
$shp.Left = EMU(3019393)
$shp.Top = EMU(1915325)
$shp.Width = EMU(5539563)
$shp.Height = EMU(584775)

# ---- Item 1 : Price ------------------------------------------------------
$shp = $s.Shapes.Item("TextBox 14")
$shp.TextFrame.TextRange.Delete()
$shp.TextFrame.TextRange.Text = "11.00"
$shp.TextFrame.TextRange.Font.Bold = $true
$shp.Left = EMU(8860212)
$shp.Top = EMU(1751561)
$shp.Width = EMU(1123760)
$shp.Height = EMU(369332)

# ---- Item 2 : Name ----------------------------------------------------
$shp = $s.Shapes.Item("TextBox 15")
$shp.TextFrame.TextRange.Delete()
$shp.TextFrame.TextRange.Text = "Fish & Chips"
$shp.TextFrame.TextRange.Font.Bold = $true
$shp.Left = EMU(3019393)
$shp.Top = EMU(2901089)
$shp.Width = EMU(5539563)
$shp.Height = EMU(369332)

# ---- Item 2 : Description ---------------------------------------------
$shp = $s.Shapes.Item("TextBox 16")
$shp.TextFrame.TextRange.Delete()
$shp.TextFrame.TextRange.Text = "beer battered fish & chips with salad & tartare sauce"
$shp.TextFrame.TextRange.Font.Size = 16
$shp.TextFrame.TextRange.Font.Color.RGB = $descGray
$shp.Left = EMU(3019393)
$shp.Top = EMU(3238522)
$shp.Width = EMU(5539563)
$shp.Height = EMU(338554)

# ---- Item 2 : Price ------------------------------------------------------
$shp = $s.Shapes.Item("TextBox 17")
$shp.TextFrame.TextRange.Delete()
$shp.TextFrame.TextRange.Text = "11.00"
$shp.TextFrame.TextRange.Font.Bold = $true

# ---- Item 3 : Name ----------------------------------------------------
$shp = $s.Shapes.Item("TextBox 18")
$shp.TextFrame.TextRange.Delete()
$shp.TextFrame.TextRange.Text = "Jalapeno Fish Tails & Chips"
$shp.TextFrame.TextRange.Font.Bold = $true

# ---- Item 3 : Description ---------------------------------------------
$shp = $s.Shapes.Item("TextBox 19")
$shp.TextFrame.TextRange.Delete()
$shp.TextFrame.TextRange.Text = "crumbed jalapeno fish tails served with chips, salad & tartare sauce"
$shp.TextFrame.TextRange.Font.Size = 16
$shp.TextFrame.TextRange.Font.Color.RGB = $descGray
$shp.Left = EMU(3019393)
$shp.Top = EMU(4401143)
$shp.Width = EMU(5539563)
$shp.Height = EMU(584775)

# ---- Item 3 : Price ------------------------------------------------------
$shp = $s.Shapes.Item("TextBox 20")
$shp.TextFrame.TextRange.Delete()
$shp.TextFrame.TextRange.Text = "11.00"
$shp.TextFrame.TextRange.Font.Bold = $true
